# Auto-generated script applying cell-level numeric updates to the Lamia_Profits workbook
# across worksheets ALC, BSM, CRP, CUL, GSM, LTW, WVR, based on the supplied OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 96
$ws.Range("I2").Value = 94.5
$ws.Range("J2").Value = 102
$ws.Range("K2").Value = 94.5
$ws.Range("L2").Value = 102
$ws.Range("M2").Value = 18.5
$ws.Range("N2").Value = -328
$ws.Range("H9").Value = 50
$ws.Range("I9").Value = 50
$ws.Range("K9").Value = 50
$ws.Range("M9").Value = 119
$ws.Range("H19").Value = 881.86664
$ws.Range("I19").Value = 520.7273
$ws.Range("K19").Value = 520.7273
$ws.Range("M19").Value = -345.7273
$ws.Range("H43").Value = 6220.875
$ws.Range("I43").Value = 4090.6667
$ws.Range("J43").Value = 7499
$ws.Range("K43").Value = 4090.6667
$ws.Range("L43").Value = 7499
$ws.Range("M43").Value = -4021.6667
$ws.Range("N43").Value = -7637
$ws.Range("H55").Value = 493
$ws.Range("I55").Value = 361.4
$ws.Range("J55").Value = 587
$ws.Range("K55").Value = 361.4
$ws.Range("L55").Value = 587
$ws.Range("M55").Value = -147.4
$ws.Range("N55").Value = -1015
$ws.Range("H58").Value = 2524.8333
$ws.Range("J58").Value = 1691.6666
$ws.Range("L58").Value = 5074.9998
$ws.Range("N58").Value = -5374.9998
$ws.Range("H70").Value = 5670.486
$ws.Range("I70").Value = 4774.1763
$ws.Range("J70").Value = 6517
$ws.Range("K70").Value = 14322.5289
$ws.Range("L70").Value = 19551
$ws.Range("M70").Value = -14052.5289
$ws.Range("N70").Value = -20091
$ws.Range("H73").Value = 5670.486
$ws.Range("I73").Value = 4774.1763
$ws.Range("J73").Value = 6517
$ws.Range("K73").Value = 14322.5289
$ws.Range("L73").Value = 19551
$ws.Range("M73").Value = -13386.5289
$ws.Range("N73").Value = -21423
$ws.Range("H86").Value = 6416
$ws.Range("I86").Value = 6261
$ws.Range("J86").Value = 6500.5454
$ws.Range("K86").Value = 6261
$ws.Range("L86").Value = 6500.5454
$ws.Range("M86").Value = -5138
$ws.Range("N86").Value = -8746.545399999999
$ws.Range("H89").Value = 6416
$ws.Range("I89").Value = 6261
$ws.Range("J89").Value = 6500.5454
$ws.Range("K89").Value = 31305
$ws.Range("L89").Value = 32502.727
$ws.Range("M89").Value = -25689
$ws.Range("N89").Value = -43734.727
$ws.Range("H116").Value = 11032.143
$ws.Range("J116").Value = 11163.667
$ws.Range("L116").Value = 11163.667
$ws.Range("N116").Value = -18047.667
$ws.Range("H125").Value = 926.6
$ws.Range("J125").Value = 1000
$ws.Range("L125").Value = 9000
$ws.Range("N125").Value = -13920
$ws.Range("H135").Value = 556.5
$ws.Range("I135").Value = 556.7368
$ws.Range("J135").Value = 555
$ws.Range("K135").Value = 5010.6312
$ws.Range("L135").Value = 4995
$ws.Range("M135").Value = -2475.6312
$ws.Range("N135").Value = -10065
$ws.Range("H138").Value = 3490.875
$ws.Range("J138").Value = 3565.7222
$ws.Range("L138").Value = 10697.1666
$ws.Range("N138").Value = -20977.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1345.25
$ws.Range("I99").Value = 1345.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1345.25
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 152.75
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1805.8334
$ws.Range("J16").Value = 1786.6666
$ws.Range("L16").Value = 1786.6666
$ws.Range("N16").Value = -2360.6666
$ws.Range("H31").Value = 42368.758
$ws.Range("I31").Value = 3393.7334
$ws.Range("J31").Value = 84127.71000000001
$ws.Range("K31").Value = 3393.7334
$ws.Range("L31").Value = 84127.71000000001
$ws.Range("M31").Value = -3098.7334
$ws.Range("N31").Value = -84717.71000000001
$ws.Range("H34").Value = 42368.758
$ws.Range("I34").Value = 3393.7334
$ws.Range("J34").Value = 84127.71000000001
$ws.Range("K34").Value = 3393.7334
$ws.Range("L34").Value = 84127.71000000001
$ws.Range("M34").Value = -3191.7334
$ws.Range("N34").Value = -84531.71000000001
$ws.Range("H52").Value = 62000
$ws.Range("J52").Value = 62000
$ws.Range("L52").Value = 62000
$ws.Range("N52").Value = -62588
$ws.Range("H99").Value = 2720.8333
$ws.Range("I99").Value = 1581.25
$ws.Range("K99").Value = 1581.25
$ws.Range("M99").Value = -83.25
$ws.Range("H103").Value = 36666.668
$ws.Range("I103").Value = 36666.668
$ws.Range("K103").Value = 36666.668
$ws.Range("M103").Value = -35494.668
$ws.Range("H107").Value = 857.57574
$ws.Range("I107").Value = 635.2083
$ws.Range("K107").Value = 635.2083
$ws.Range("M107").Value = 1284.7917
$ws.Range("H113").Value = 1805.8334
$ws.Range("J113").Value = 1786.6666
$ws.Range("L113").Value = 1786.6666
$ws.Range("N113").Value = -6126.6666
$ws.Range("H122").Value = 7405.467
$ws.Range("J122").Value = 13873
$ws.Range("L122").Value = 41619
$ws.Range("N122").Value = -46519
$ws.Range("H126").Value = 2720.8333
$ws.Range("I126").Value = 1581.25
$ws.Range("K126").Value = 4743.75
$ws.Range("M126").Value = -2273.75
$ws.Range("H132").Value = 5275.6924
$ws.Range("I132").Value = 5111.5625
$ws.Range("J132").Value = 5739.1177
$ws.Range("K132").Value = 15334.6875
$ws.Range("L132").Value = 17217.3531
$ws.Range("M132").Value = -12804.6875
$ws.Range("N132").Value = -22277.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 469.625
$ws.Range("J23").Value = 610.25
$ws.Range("L23").Value = 1830.75
$ws.Range("N23").Value = -2300.75
$ws.Range("H42").Value = 10450
$ws.Range("J42").Value = 10450
$ws.Range("L42").Value = 31350
$ws.Range("N42").Value = -32418
$ws.Range("H47").Value = 33901
$ws.Range("I47").Value = 33901
$ws.Range("K47").Value = 101703
$ws.Range("M47").Value = -101272
$ws.Range("H56").Value = 9646.362999999999
$ws.Range("I56").Value = 9646.362999999999
$ws.Range("K56").Value = 9646.362999999999
$ws.Range("M56").Value = -9116.362999999999
$ws.Range("H69").Value = 10002.8
$ws.Range("J69").Value = 10002.8
$ws.Range("L69").Value = 30008.4
$ws.Range("N69").Value = -31630.4
$ws.Range("H72").Value = 10002.8
$ws.Range("J72").Value = 10002.8
$ws.Range("L72").Value = 90025.2
$ws.Range("N72").Value = -98137.2
$ws.Range("H87").Value = 24800.867
$ws.Range("I87").Value = 20506.5
$ws.Range("J87").Value = 25461.54
$ws.Range("K87").Value = 61519.5
$ws.Range("L87").Value = 76384.62
$ws.Range("M87").Value = -60271.5
$ws.Range("N87").Value = -78880.62
$ws.Range("H90").Value = 24800.867
$ws.Range("I90").Value = 20506.5
$ws.Range("J90").Value = 25461.54
$ws.Range("K90").Value = 184558.5
$ws.Range("L90").Value = 229153.86
$ws.Range("M90").Value = -178318.5
$ws.Range("N90").Value = -241633.86

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 919.5405
$ws.Range("J97").Value = 1017.7692
$ws.Range("L97").Value = 1017.7692
$ws.Range("N97").Value = -2009.7692
$ws.Range("H107").Value = 1969
$ws.Range("I107").Value = 821.7143
$ws.Range("K107").Value = 821.7143
$ws.Range("M107").Value = 1098.2857
$ws.Range("H132").Value = 44951.668
$ws.Range("I132").Value = 50773.13
$ws.Range("J132").Value = 11478.25
$ws.Range("K132").Value = 152319.39
$ws.Range("L132").Value = 34434.75
$ws.Range("M132").Value = -149789.39
$ws.Range("N132").Value = -39494.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8236.388999999999
$ws.Range("I7").Value = 2848.4614
$ws.Range("J7").Value = 22245
$ws.Range("K7").Value = 2848.4614
$ws.Range("L7").Value = 22245
$ws.Range("M7").Value = -2736.4614
$ws.Range("N7").Value = -22469
$ws.Range("H46").Value = 5300.25
$ws.Range("I46").Value = 3700
$ws.Range("J46").Value = 5833.6665
$ws.Range("K46").Value = 3700
$ws.Range("L46").Value = 5833.6665
$ws.Range("M46").Value = -3512
$ws.Range("N46").Value = -6209.6665
$ws.Range("H55").Value = 1667195
$ws.Range("I55").Value = 2500451.8
$ws.Range("K55").Value = 2500451.8
$ws.Range("M55").Value = -2500278.8
$ws.Range("H122").Value = 5113.552
$ws.Range("I122").Value = 4087.08
$ws.Range("K122").Value = 12261.24
$ws.Range("M122").Value = -9811.24
$ws.Range("H126").Value = 8236.388999999999
$ws.Range("I126").Value = 2848.4614
$ws.Range("J126").Value = 22245
$ws.Range("K126").Value = 8545.3842
$ws.Range("L126").Value = 66735
$ws.Range("M126").Value = -6075.3842
$ws.Range("N126").Value = -71675
$ws.Range("H132").Value = 2769.0833
$ws.Range("I132").Value = 1498.6296
$ws.Range("J132").Value = 6580.4443
$ws.Range("K132").Value = 4495.8888
$ws.Range("L132").Value = 19741.3329
$ws.Range("M132").Value = -1965.8888
$ws.Range("N132").Value = -24801.3329
$ws.Range("H136").Value = 4108.478
$ws.Range("J136").Value = 10834.833
$ws.Range("L136").Value = 32504.499
$ws.Range("N136").Value = -37604.499
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 14550
$ws.Range("H107").Value = 415.125
$ws.Range("I107").Value = 274.6
$ws.Range("K107").Value = 823.8000000000001
$ws.Range("M107").Value = 1096.2
$ws.Range("H126").Value = 5820.3687
$ws.Range("I126").Value = 5477.0557
$ws.Range("J126").Value = 12000
$ws.Range("K126").Value = 16431.1671
$ws.Range("L126").Value = 36000
$ws.Range("M126").Value = -13961.1671
$ws.Range("N126").Value = -40940

